# Updated cryptos list on Tue Feb 27 14:01:16 UTC 2024 with GitHub Actions
# Applies updated Price (D) and Volume(1h) (E) values, and fixes the
# Kaspa/RenderToken row ordering (rows 29-30: Coin/Link/Price/Volume).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.093.55"
$ws.Range("E2").Value = "'  +11.42%  "
$ws.Range("D3").Value = "'3.268.76"
$ws.Range("E3").Value = "'  +6.78%  "
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("D5").Value = "'397.93"
$ws.Range("E5").Value = "'  +1.00%  "
$ws.Range("D6").Value = "'110.26"
$ws.Range("E6").Value = "'  +8.56%  "
$ws.Range("D7").Value = "'0.562"
$ws.Range("E7").Value = "'  +5.55%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "'  -0.04%  "
$ws.Range("D9").Value = "'0.627"
$ws.Range("E9").Value = "'  +7.50%  "
$ws.Range("D10").Value = "'39.30"
$ws.Range("E10").Value = "'  +6.27%  "
$ws.Range("D11").Value = "'0.0981"
$ws.Range("E11").Value = "'  +15.53%  "
$ws.Range("E12").Value = "'  +2.50%  "
$ws.Range("D13").Value = "'3.786.02"
$ws.Range("E13").Value = "'  +6.69%  "
$ws.Range("D14").Value = "'8.19"
$ws.Range("E14").Value = "'  +6.81%  "
$ws.Range("D15").Value = "'19.11"
$ws.Range("E15").Value = "'  +3.45%  "
$ws.Range("D16").Value = "'3.267.13"
$ws.Range("E16").Value = "'  +6.77%  "
$ws.Range("E17").Value = "'  +2.69%  "
$ws.Range("D18").Value = "'10.83"
$ws.Range("E18").Value = "'  +3.22%  "
$ws.Range("D19").Value = "'56.919.81"
$ws.Range("E19").Value = "'  +11.06%  "
$ws.Range("E20").Value = "'  +5.66%  "
$ws.Range("D21").Value = "'0.0000107"
$ws.Range("E21").Value = "'  +11.82%  "
$ws.Range("D22").Value = "'12.96"
$ws.Range("E22").Value = "'  +5.79%  "
$ws.Range("D23").Value = "'307.71"
$ws.Range("E23").Value = "'  +16.53%  "
$ws.Range("D24").Value = "'75.49"
$ws.Range("E24").Value = "'  +8.07%  "
$ws.Range("E25").Value = "'  -0.30%  "
$ws.Range("D26").Value = "'28.27"
$ws.Range("E26").Value = "'  +5.35%  "
$ws.Range("D27").Value = "'7.92"
$ws.Range("E27").Value = "'  +1.10%  "
$ws.Range("D28").Value = "'4.39"
$ws.Range("E28").Value = "'  +5.25%  "
$ws.Range("B29").Value = "'Kaspa"
$ws.Range("C29").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.170"
$ws.Range("E29").Value = "'  +4.59%  "
$ws.Range("B30").Value = "'RenderToken"
$ws.Range("C30").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "'7.26"
$ws.Range("E30").Value = "'  +1.76%  "
$ws.Range("E31").Value = "'  +0.16%  "
$ws.Range("E32").Value = "'  +4.64%  "
$ws.Range("D33").Value = "'11.06"
$ws.Range("E33").Value = "'  +4.45%  "
$ws.Range("D34").Value = "'37.62"
$ws.Range("E34").Value = "'  +4.82%  "
$ws.Range("D35").Value = "'0.0481"
$ws.Range("E35").Value = "'  -3.74%  "
$ws.Range("D36").Value = "'2.14"
$ws.Range("E36").Value = "'  +3.19%  "
$ws.Range("D37").Value = "'51.56"
$ws.Range("E37").Value = "'  +3.06%  "
$ws.Range("D38").Value = "'3.15"
$ws.Range("E38").Value = "'  +23.97%  "
$ws.Range("D39").Value = "'3.55"
$ws.Range("E39").Value = "'  +8.04%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "'  -0.13%  "
$ws.Range("D41").Value = "'135.18"
$ws.Range("E41").Value = "'  +5.25%  "
$ws.Range("E42").Value = "'  +5.47%  "
$ws.Range("E43").Value = "'  +5.26%  "
$ws.Range("E44").Value = "'  +4.82%  "
$ws.Range("D45").Value = "'3.95"
$ws.Range("E45").Value = "'  +1.04%  "
$ws.Range("E46").Value = "'  -2.66%  "
$ws.Range("D47").Value = "'22.12"
$ws.Range("E47").Value = "'  +2.49%  "
$ws.Range("D48").Value = "'2.153.75"
$ws.Range("E48").Value = "'  +4.38%  "
$ws.Range("E49").Value = "'  +2.52%  "
$ws.Range("D50").Value = "'2.04"
$ws.Range("E50").Value = "'  +42.12%  "
$ws.Range("E51").Value = "'  -4.67%  "
